$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 129 (shifts existing rows 129:180 down to 130:181)
$ws.Rows.Item(129).Insert()

# Populate the newly inserted row 129 with the new data record
$ws.Cells.Item(129, 1).Value = 10
$ws.Cells.Item(129, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(129, 3).Value = "La Araucanía"
$ws.Cells.Item(129, 4).Value = 44523
$ws.Cells.Item(129, 5).Value = 9
$ws.Cells.Item(129, 6).Value = 100112039
$ws.Cells.Item(129, 7).Value = "Ciboulette"
$ws.Cells.Item(129, 8).Value = "Sin especificar"
$ws.Cells.Item(129, 9).Value = "Primera"
$ws.Cells.Item(129, 10).Value = 30
$ws.Cells.Item(129, 11).Value = 5000
$ws.Cells.Item(129, 12).Value = 5000
$ws.Cells.Item(129, 13).Value = 5000
$ws.Cells.Item(129, 14).Value = "`$/docena de atados"
$ws.Cells.Item(129, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(129, 16).Value = 1667
$ws.Cells.Item(129, 17).Value = 3
$ws.Cells.Item(129, 18).Value = "Hortaliza"
